$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "30.269.21"
    "E2" = "  +0.09%  "
    "D3" = "2.001.87"
    "E3" = "  +5.67%  "
    "D4" = "1.002"
    "E4" = "  +0.03%  "
    "D5" = "323.54"
    "E5" = "  +0.74%  "
    "D7" = "0.5100"
    "E7" = "  +0.59%  "
    "D8" = "0.4148"
    "E8" = "  +2.61%  "
    "D9" = "0.08706"
    "E9" = "  +4.88%  "
    "E10" = "  +2.19%  "
    "D11" = "42.74"
    "E11" = "  +1.17%  "
    "D12" = "24.87"
    "E12" = "  +2.49%  "
    "D13" = "2.000.70"
    "E13" = "  +5.44%  "
    "D14" = "6.534"
    "E14" = "  +1.95%  "
    "D15" = "7.444"
    "E15" = "  +1.51%  "
    "D16" = "1.003"
    "E16" = "  +0.10%  "
    "D17" = "94.15"
    "E17" = "  +1.26%  "
    "E18" = "  +1.49%  "
    "D19" = "0.06545"
    "E19" = "  +1.24%  "
    "D20" = "18.94"
    "E21" = "  -0.06%  "
    "D22" = "6.136"
    "E22" = "  +3.57%  "
    "D23" = "30.330.64"
    "E23" = "  +0.26%  "
    "D24" = "11.72"
    "E24" = "  +3.71%  "
    "D25" = "2.208"
    "E25" = "  +1.08%  "
    "D26" = "2.233.92"
    "E26" = "  +5.79%  "
    "E27" = "  +4.24%  "
    "D28" = "163.76"
    "E28" = "  +1.68%  "
    "D29" = "2.406"
    "E29" = "  +5.63%  "
    "D30" = "131.26"
    "E30" = "  +1.75%  "
    "D31" = "1.140"
    "E31" = "  +2.64%  "
    "D32" = "0.1051"
    "E32" = "  +0.70%  "
    "D33" = "6.069"
    "E33" = "  +0.89%  "
    "D34" = "3.832"
    "E34" = "  +3.08%  "
    "D35" = "1.339"
    "E35" = "  +12.63%  "
    "D36" = "0.02503"
    "E36" = "  +2.01%  "
    "B37" = "Hedera"
    "C37" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D37" = "0.06593"
    "E37" = "  +2.14%  "
    "B38" = "InternetComputer(DFINITY)"
    "C38" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D38" = "5.388"
    "E38" = "  +0.92%  "
    "D39" = "12.28"
    "E39" = "  +7.67%  "
    "D40" = "0.2200"
    "E40" = "  +1.97%  "
    "D41" = "8.929"
    "E41" = "  +3.55%  "
    "D42" = "0.6632"
    "E42" = "  +3.38%  "
    "D43" = "1.228"
    "E43" = "  +1.13%  "
    "E44" = "  +2.76%  "
    "D45" = "0.6153"
    "E45" = "  +2.83%  "
    "D46" = "2.200"
    "E46" = "  +2.44%  "
    "D47" = "3.659"
    "E47" = "  +0.43%  "
    "D48" = "1.261"
    "E48" = "  +3.74%  "
    "D49" = "124.63"
    "E49" = "  +0.65%  "
    "D50" = "79.97"
    "E50" = "  +1.46%  "
    "D51" = "0.06891"
    "E51" = "  +1.82%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
